# Updated cryptos list on Wed May 31 23:20:50 UTC 2023 with GitHub Actions
#
# Refreshes the scraped "Price" (column D) and "Volume(1h)" (column E)
# figures for each coin row on the single worksheet, and re-syncs rows
# 46/47 (Decentraland / PaxDollar) whose relative order changed in the
# latest coinranking.com pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks like a number
# (e.g. "1.000", "306.42"), so it stays a string cell exactly like the
# other inline-string cells on this sheet instead of becoming a numeric
# cell.
function Set-TextValue($CellRef, $Text) {
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "27.167.57"
$ws.Range("E2").Value = "  -1.96%  "
# Row 3
$ws.Range("D3").Value = "1.869.23"
$ws.Range("E3").Value = "  -1.65%  "
# Row 4
Set-TextValue "D4" "1.000"
# Row 5
Set-TextValue "D5" "306.42"
$ws.Range("E5").Value = "  -1.83%  "
# Row 6
Set-TextValue "D6" "0.9997"
$ws.Range("E6").Value = "  -0.12%  "
# Row 7
Set-TextValue "D7" "0.5189"
$ws.Range("E7").Value = "  -0.33%  "
# Row 8
Set-TextValue "D8" "0.3743"
$ws.Range("E8").Value = "  -1.10%  "
# Row 9
$ws.Range("E9").Value = "  -1.22%  "
# Row 10
Set-TextValue "D10" "0.8921"
$ws.Range("E10").Value = "  -1.08%  "
# Row 11
$ws.Range("E11").Value = "  -1.98%  "
# Row 12
$ws.Range("D12").Value = "1.875.15"
$ws.Range("E12").Value = "  -1.25%  "
# Row 13
Set-TextValue "D13" "0.07520"
$ws.Range("E13").Value = "  -1.62%  "
# Row 14
Set-TextValue "D14" "5.307"
$ws.Range("E14").Value = "  -2.56%  "
# Row 15
Set-TextValue "D15" "90.05"
$ws.Range("E15").Value = "  -2.23%  "
# Row 16
Set-TextValue "D16" "1.000"
$ws.Range("E16").Value = "  -0.15%  "
# Row 17
$ws.Range("E17").Value = "  -2.38%  "
# Row 18
Set-TextValue "D18" "14.12"
$ws.Range("E18").Value = "  -2.32%  "
# Row 19
Set-TextValue "D19" "0.9996"
$ws.Range("E19").Value = "  -0.06%  "
# Row 20
$ws.Range("D20").Value = "27.194.36"
$ws.Range("E20").Value = "  -2.00%  "
# Row 21
Set-TextValue "D21" "5.006"
# Row 22
$ws.Range("D22").Value = "2.114.28"
$ws.Range("E22").Value = "  -0.71%  "
# Row 23
Set-TextValue "D23" "10.46"
$ws.Range("E23").Value = "  -3.49%  "
# Row 24
$ws.Range("E24").Value = "  -2.11%  "
# Row 25
Set-TextValue "D25" "1.835"
$ws.Range("E25").Value = "  -1.70%  "
# Row 26
$ws.Range("E26").Value = "  -4.55%  "
# Row 27
Set-TextValue "D27" "17.97"
$ws.Range("E27").Value = "  -1.80%  "
# Row 28
Set-TextValue "D28" "2.083"
$ws.Range("E28").Value = "  -3.39%  "
# Row 29
$ws.Range("E29").Value = "  -1.36%  "
# Row 30
Set-TextValue "D30" "4.663"
$ws.Range("E30").Value = "  -3.65%  "
# Row 31
Set-TextValue "D31" "4.690"
$ws.Range("E31").Value = "  -3.19%  "
# Row 32
Set-TextValue "D32" "0.09259"
$ws.Range("E32").Value = "  +2.01%  "
# Row 33
Set-TextValue "D33" "0.05133"
$ws.Range("E33").Value = "  -3.05%  "
# Row 34
Set-TextValue "D34" "3.084"
$ws.Range("E34").Value = "  -3.32%  "
# Row 35
Set-TextValue "D35" "1.160"
$ws.Range("E35").Value = "  -5.39%  "
# Row 36
Set-TextValue "D36" "0.7276"
$ws.Range("E36").Value = "  -6.58%  "
# Row 37
Set-TextValue "D37" "0.02034"
$ws.Range("E37").Value = "  -2.59%  "
# Row 38
Set-TextValue "D38" "3.118"
$ws.Range("E38").Value = "  +1.48%  "
# Row 39
Set-TextValue "D39" "2.509"
$ws.Range("E39").Value = "  -2.82%  "
# Row 40
$ws.Range("E40").Value = "  -1.72%  "
# Row 41
Set-TextValue "D41" "0.5313"
$ws.Range("E41").Value = "  -4.27%  "
# Row 42
Set-TextValue "D42" "6.532"
$ws.Range("E42").Value = "  -2.82%  "
# Row 43
Set-TextValue "D43" "116.24"
$ws.Range("E43").Value = "  +0.25%  "
# Row 44
Set-TextValue "D44" "8.322"
$ws.Range("E44").Value = "  -2.06%  "
# Row 45
Set-TextValue "D45" "0.1474"
$ws.Range("E45").Value = "  -2.64%  "
# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.4631"
$ws.Range("E46").Value = "  -3.75%  "
# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D47" "0.9992"
$ws.Range("E47").Value = "  -0.15%  "
# Row 48
Set-TextValue "D48" "10.02"
$ws.Range("E48").Value = "  -3.95%  "
# Row 49
Set-TextValue "D49" "1.564"
$ws.Range("E49").Value = "  -2.90%  "
# Row 50
Set-TextValue "D50" "36.74"
$ws.Range("E50").Value = "  -0.66%  "
# Row 51
Set-TextValue "D51" "63.67"
$ws.Range("E51").Value = "  -4.55%  "
